$d = $word.ActiveDocument

# --- 1) Add an extra manual line break after "{{ formType }}" -------------
$rng = $d.Content
$rng.Find.Execute("{{ formType }}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{{ formType }}^l", 2) | Out-Null

# --- 2) Add the modifiedAt / generatedAt placeholders after createdAt -----
$rng = $d.Content
$rng.Find.Execute("{{ createdAt }}", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{{ createdAt }}^l{{ modifiedAt }}^l{{ generatedAt }}", 2) | Out-Null

# --- 3) Shrink column 1 / widen column 2 of the first table (all rows) ----
# 3871 twips -> 3869 twips  (193.55pt -> 193.45pt)
# 272  twips -> 274  twips  (13.6pt  -> 13.7pt)
$t = $d.Tables.Item(1)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $t.Cell($r, 1).Width = 193.45
    $t.Cell($r, 2).Width = 13.7
}
